$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Formula = "'" + $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '68.049.40'
Set-TextCell 'E2' '  -0.47%  '

Set-TextCell 'D3' '3.638.24'
Set-TextCell 'E3' '  -1.09%  '

Set-TextCell 'D4' '1.00'
Set-TextCell 'E4' '  +0.26%  '

Set-TextCell 'D5' '587.57'
Set-TextCell 'E5' '  -1.89%  '

Set-TextCell 'D6' '195.49'
Set-TextCell 'E6' '  +0.21%  '

Set-TextCell 'D7' '3.630.88'
Set-TextCell 'E7' '  -1.04%  '

Set-TextCell 'E8' '  -0.46%  '

Set-TextCell 'D9' '1.00'
Set-TextCell 'E9' '  +0.06%  '

Set-TextCell 'E10' '  -2.63%  '

Set-TextCell 'D11' '0.154'
Set-TextCell 'E11' '  +0.03%  '

Set-TextCell 'D12' '55.52'
Set-TextCell 'E12' '  -4.09%  '

Set-TextCell 'D13' '0.0000292'
Set-TextCell 'E13' '  +6.33%  '

Set-TextCell 'D14' '10.05'
Set-TextCell 'E14' '  -2.29%  '

Set-TextCell 'D15' '4.196.40'
Set-TextCell 'E15' '  -1.33%  '

Set-TextCell 'D16' '3.629.90'
Set-TextCell 'E16' '  -1.04%  '

Set-TextCell 'E17' '  -0.64%  '

Set-TextCell 'E18' '  -0.85%  '

Set-TextCell 'B19' 'WrappedBTC'
Set-TextCell 'C19' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 'D19' '67.964.47'
Set-TextCell 'E19' '  -0.25%  '

Set-TextCell 'B20' 'Chainlink'
Set-TextCell 'C20' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 'D20' '18.54'
Set-TextCell 'E20' '  -2.39%  '

Set-TextCell 'D21' '1.08'
Set-TextCell 'E21' '  -2.75%  '

Set-TextCell 'D22' '405.41'
Set-TextCell 'E22' '  +0.64%  '

Set-TextCell 'D23' '13.37'
Set-TextCell 'E23' '  +18.85%  '

Set-TextCell 'E24' '  -3.50%  '

Set-TextCell 'D25' '86.36'
Set-TextCell 'E25' '  -1.97%  '

Set-TextCell 'D26' '2.95'
Set-TextCell 'E26' '  -0.90%  '

Set-TextCell 'D27' '12.65'
Set-TextCell 'E27' '  -0.82%  '

Set-TextCell 'D28' '3.91'
Set-TextCell 'E28' '  +5.34%  '

Set-TextCell 'E29' '  +0.51%  '

Set-TextCell 'D30' '8.16'
Set-TextCell 'E30' '  +11.18%  '

Set-TextCell 'D31' '9.26'
Set-TextCell 'E31' '  -1.70%  '

Set-TextCell 'D32' '31.63'

Set-TextCell 'D33' '675.20'
Set-TextCell 'E33' '  +10.54%  '

Set-TextCell 'D34' '12.29'
Set-TextCell 'E34' '  -0.80%  '

Set-TextCell 'E35' '  +1.32%  '

Set-TextCell 'D36' '64.92'
Set-TextCell 'E36' '  -2.10%  '

Set-TextCell 'D37' '43.01'
Set-TextCell 'E37' '  -5.29%  '

Set-TextCell 'D38' '0.423'
Set-TextCell 'E38' '  +6.19%  '

Set-TextCell 'D39' '0.0₃0814'
Set-TextCell 'E39' '  +4.62%  '

Set-TextCell 'E40' '  +0.13%  '

Set-TextCell 'D41' '2.94'
Set-TextCell 'E41' '  +15.26%  '

Set-TextCell 'D42' '3.188.99'
Set-TextCell 'E42' '  +14.06%  '

Set-TextCell 'D43' '3.10'
Set-TextCell 'E43' '  +6.57%  '

Set-TextCell 'E44' '  -1.68%  '

Set-TextCell 'E45' '  -0.13%  '

Set-TextCell 'D46' '0.0421'
Set-TextCell 'E46' '  -1.98%  '

Set-TextCell 'D47' '8.90'
Set-TextCell 'E47' '  -0.01%  '

Set-TextCell 'E48' '  -2.85%  '

Set-TextCell 'B49' 'Monero'
Set-TextCell 'C49' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D49' '143.36'
Set-TextCell 'E49' '  -0.46%  '

Set-TextCell 'B50' 'ApeXProtocol'
Set-TextCell 'C50' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 'D50' '3.09'
Set-TextCell 'E50' '  -2.33%  '

Set-TextCell 'B51' 'dogwifhat'
Set-TextCell 'C51' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D51' '2.55'
Set-TextCell 'E51' '  +0.52%  '
